$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "dry_15:16"
$ws.Range("A3").Value = "dry_16:17"
$ws.Range("A4").Value = "dry_17:18"
$ws.Range("A5").Value = "dry_18:19"
$ws.Range("A6").Value = "dry_19:20"
$ws.Range("A7").Value = "dry_20:21"
$ws.Range("A8").Value = "dry_21:22"
$ws.Range("A9").Value = "dry_22:23"
$ws.Range("A10").Value = "dry_23:24"
$ws.Range("A11").Value = "rainy_15:16"
$ws.Range("A12").Value = "rainy_16:17"
$ws.Range("A13").Value = "rainy_17:18"
$ws.Range("A14").Value = "rainy_18:19"
$ws.Range("A15").Value = "rainy_19:20"
$ws.Range("A16").Value = "rainy_20:21"
$ws.Range("A17").Value = "rainy_21:22"
$ws.Range("A18").Value = "rainy_22:23"
$ws.Range("A19").Value = "rainy_23:24"
